# Apply the edits described by the commit to the workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet1: fill in the "Rest API" column (D) for rows 9-14 ---------------
$ws1 = $wb.Worksheets.Item("Sheet1")

# Note: values are entered in an order that makes the resulting shared-string
# table match the target layout (snapshots, sites, instanceView, then the
# "Refer to..." note).
$ws1.Range("D13").Value = "GET https://management.azure.com/subscriptions/{subscriptionId}/resourceGroups/{resourceGroupName}/providers/Microsoft.Compute/snapshots/{snapshotName}?api-version=2023-01-01"
$ws1.Range("D14").Value = "GET https://management.azure.com/subscriptions/{subscriptionId}/resourceGroups/{resourceGroupName}/providers/Microsoft.Web/sites?api-version=2023-01-01"
$ws1.Range("D9").Value  = "GET https://management.azure.com/subscriptions/{subscriptionId}/resourceGroups/{resourceGroupName}/providers/Microsoft.Compute/virtualMachines/{vmName}/instanceView?api-version=2023-01-01"
$ws1.Range("D10").Value = "GET https://management.azure.com/subscriptions/{subscriptionId}/resourceGroups/{resourceGroupName}/providers/Microsoft.Compute/virtualMachines/{vmName}/instanceView?api-version=2023-01-01"
$ws1.Range("D11").Value = "Refer to the Azure API documentation"
$ws1.Range("D12").Value = "Refer to the Azure API documentation"

# --- Sheet1: move the current selection from D12 up to D5 -------------------
$ws1.Range("D5").Select()

# --- Sheet2: shrink the row heights for rows 7 and 8 -------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Rows.Item(7).RowHeight = 90
$ws2.Rows.Item(8).RowHeight = 60
